$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. Widen the second column (998 twips -> 2194 twips == 49.9pt -> 109.7pt)
$t.Columns.Item(2).Width = 109.7

# 2. The "4" row currently carries the hidden "_GoBack" bookmark as a direct
#    child of <w:tr> (not inside a <w:tc>). Word's object model hides
#    "_GoBack" from the Bookmarks collection, so it can't be relocated
#    in place. Delete the whole row (this removes the old bookmark too)
#    and rebuild its contents as a fresh row.
$row4 = $t.Rows.Item(9)
$row4.Delete()

$refRow = $t.Rows.Item(9)
$newRow4 = $t.Rows.Add($refRow)
$newRow4.Cells.Item(1).Range.Text = "4"
$newRow4.Cells.Item(2).Range.Text = "Low/High Gear Toggle"

# 3. Insert a brand new row (the merged Low/High Gear command becomes
#    button "5"), pushing the old numbered rows down by one.
$refRow2 = $t.Rows.Item(10)
$newRow5 = $t.Rows.Add($refRow2)
$newRow5.Cells.Item(1).Range.Text = "5"

# 4. Renumber the remaining rows: 5->6, 6->7, 7->8, 8->9, 9->10
for ($i = 11; $i -le 15; $i++) {
    $cell = $t.Rows.Item($i).Cells.Item(1)
    $txt = $cell.Range.Text
    $clean = $txt.Substring(0, $txt.Length - 2)
    $newVal = [int]$clean + 1
    $cell.Range.Text = [string]$newVal
}

# 5. Re-create the "_GoBack" bookmark, now collapsed inside the empty
#    paragraph of the last row's second cell.
$lastRow = $t.Rows.Item($t.Rows.Count)
$lastCell = $lastRow.Cells.Item(2)
$rng = $lastCell.Range
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null
